$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.698.49"
$ws.Range("E2").Value = "  +0.70%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.124.63"
$ws.Range("E3").Value = "  +1.11%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +1.02%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.97"
$ws.Range("E5").Value = "  +2.11%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.97%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5281"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4576"
$ws.Range("E8").Value = "  +2.00%  "

# Row 9 - OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.81"
$ws.Range("E9").Value = "  +1.87%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09123"
$ws.Range("E10").Value = "  +2.43%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  +2.06%  "

# Row 12 - Solana
$ws.Range("E12").Value = "  +1.01%  "

# Row 13 - WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.123.63"
$ws.Range("E13").Value = "  +1.33%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.872"
$ws.Range("E14").Value = "  +2.08%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.132"
$ws.Range("E15").Value = "  +5.42%  "

# Row 16 - ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001179"
$ws.Range("E16").Value = "  +4.88%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.40"

# Row 18 - BinanceUSD
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.014"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19 - TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06708"
$ws.Range("E19").Value = "  +1.35%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.59"
$ws.Range("E20").Value = "  +1.98%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +1.10%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.501"
$ws.Range("E22").Value = "  +3.45%  "

# Row 23 - WrappedBTC
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.772.93"
$ws.Range("E23").Value = "  +0.79%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  +5.58%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.360"
$ws.Range("E25").Value = "  +1.82%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.374.22"
$ws.Range("E26").Value = "  +1.45%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +1.17%  "

# Row 28 - Monero
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.46"
$ws.Range("E28").Value = "  +0.93%  "

# Row 29 - LidoDAOToken
$ws.Range("E29").Value = "  -0.47%  "

# Row 30 - BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.91"
$ws.Range("E30").Value = "  +2.75%  "

# Row 31 - ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.216"
$ws.Range("E31").Value = "  +1.78%  "

# Row 32 - Stellar
$ws.Range("E32").Value = "  +0.71%  "

# Row 33 - ARBITRUM
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.662"
$ws.Range("E33").Value = "  +0.14%  "

# Row 34 - Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.395"
$ws.Range("E34").Value = "  +4.01%  "

# Row 35 - HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.949"
$ws.Range("E35").Value = "  +0.22%  "

# Row 36 - FraxShare
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.63"
$ws.Range("E36").Value = "  +1.71%  "

# Rows 37 & 38 swap: VeChain and InternetComputer(DFINITY) exchange order
# Row 37 becomes InternetComputer(DFINITY)
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.935"
$ws.Range("E37").Value = "  +8.28%  "

# Row 38 becomes VeChain
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02685"
$ws.Range("E38").Value = "  +4.40%  "

# Row 39 - Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06900"
$ws.Range("E39").Value = "  +1.94%  "

# Row 40 - Algorand
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2335"
$ws.Range("E40").Value = "  +3.14%  "

# Row 41 - Aptos
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.67"
$ws.Range("E41").Value = "  -0.38%  "

# Row 42 - TheSandbox
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6938"
$ws.Range("E42").Value = "  +0.34%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.267"
$ws.Range("E43").Value = "  +0.85%  "

# Row 44 - EnergySwap
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.18"
$ws.Range("E44").Value = "  +8.46%  "

# Row 45 - Decentraland
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6501"
$ws.Range("E45").Value = "  +2.19%  "

# Row 46 - NEARProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.317"
$ws.Range("E46").Value = "  +1.14%  "

# Row 47 - BabyDogeCoin
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000375"
$ws.Range("E47").Value = "  +17.39%  "

# Row 48 - PancakeSwap
$ws.Range("E48").Value = "  +1.86%  "

# Row 49 - EOS
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.263"
$ws.Range("E49").Value = "  +1.64%  "

# Row 50 - Aave
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.99"
$ws.Range("E50").Value = "  +1.47%  "

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07314"
$ws.Range("E51").Value = "  +3.72%  "

